# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The worker/period detail rows (16-27) are re-grouped: instead of
# alternating between the two workers period-by-period, each worker now
# gets a contiguous block of rows (periods listed 2107 -> 2102,
# descending), starting with VICKY GONZALEZ BOLAÑOS followed by
# MILTON JOSE BOLAÑO BOLAÑO. The per-person/per-period mora and salary
# values themselves are unchanged, only the row order changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 16-27: DocType, DocNumber, Name, Periodo, ValorMora, SalarioBasico
$rows = @(
    @{ Row = 16; Doc = "1048609333"; Name = "VICKY GONZALEZ BOLAÑOS";    Periodo = "2107"; Mora = 32000; Salario = 1160000 },
    @{ Row = 17; Doc = "1048609333"; Name = "VICKY GONZALEZ BOLAÑOS";    Periodo = "2106"; Mora = 38400; Salario = 1160000 },
    @{ Row = 18; Doc = "1048609333"; Name = "VICKY GONZALEZ BOLAÑOS";    Periodo = "2105"; Mora = 38400; Salario = 1160000 },
    @{ Row = 19; Doc = "1048609333"; Name = "VICKY GONZALEZ BOLAÑOS";    Periodo = "2104"; Mora = 38400; Salario = 1160000 },
    @{ Row = 20; Doc = "1048609333"; Name = "VICKY GONZALEZ BOLAÑOS";    Periodo = "2103"; Mora = 38400; Salario = 1160000 },
    @{ Row = 21; Doc = "1048609333"; Name = "VICKY GONZALEZ BOLAÑOS";    Periodo = "2102"; Mora = 38400; Salario = 1160000 },
    @{ Row = 22; Doc = "9186461";    Name = "MILTON JOSE BOLAÑO BOLAÑO"; Periodo = "2107"; Mora = 33333; Salario = 1000000 },
    @{ Row = 23; Doc = "9186461";    Name = "MILTON JOSE BOLAÑO BOLAÑO"; Periodo = "2106"; Mora = 40000; Salario = 1000000 },
    @{ Row = 24; Doc = "9186461";    Name = "MILTON JOSE BOLAÑO BOLAÑO"; Periodo = "2105"; Mora = 40000; Salario = 1000000 },
    @{ Row = 25; Doc = "9186461";    Name = "MILTON JOSE BOLAÑO BOLAÑO"; Periodo = "2104"; Mora = 40000; Salario = 1000000 },
    @{ Row = 26; Doc = "9186461";    Name = "MILTON JOSE BOLAÑO BOLAÑO"; Periodo = "2103"; Mora = 40000; Salario = 1000000 },
    @{ Row = 27; Doc = "9186461";    Name = "MILTON JOSE BOLAÑO BOLAÑO"; Periodo = "2102"; Mora = 40000; Salario = 1000000 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = "CC"
    $ws.Range("C$n").Value = $r.Doc
    $ws.Range("D$n").Value = $r.Name
    $ws.Range("E$n").Value = $r.Periodo
    $ws.Range("F$n").Value = $r.Mora
    $ws.Range("G$n").Value = $r.Salario
}
